# categorical_summary.xlsx - apply formatting touch-up:
#  - set explicit custom widths for columns A, B, C
#  - leave the A1:B81 range selected (matches the saved selection state)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 43.5
$ws.Columns.Item(2).ColumnWidth = 63.66666666666667
$ws.Columns.Item(3).ColumnWidth = 69.83333333333334

$ws.Range("A1:B81").Select()
